# This script applies a rotation/exchange of species-record data among a
# handful of row groups on the "Artfynd" sheet: {2,3,4}, {5,6}, {9,10},
# {11,12}, {13,14} and {18,19}. Within each group the identifying columns
# (A Id, B Taxonsorteringsordning, E TaxonId, F Artnamn, G Vetenskapligt
# namn, H Auktor, Q Ost, R Nord, AC Publik kommentar) are swapped between
# rows while every other column (coordinates-independent metadata such as
# Rödlistade, Lokalnamn, datum, observatör, etc.) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 131064784
$ws.Range("B2").Value = 91828
$ws.Range("E2").Value = 5432
$ws.Range("F2").Value = 'Granticka'
$ws.Range("G2").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H2").ClearContents()
$ws.Range("Q2").Value = 442100
$ws.Range("R2").Value = 7039221
$ws.Range("AC2").ClearContents()

# Row 3
$ws.Range("A3").Value = 131064775
$ws.Range("B3").Value = 57884
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = 'Tretåig hackspett'
$ws.Range("G3").Value = 'Picoides tridactylus'
$ws.Range("H3").Value = '(Linnaeus, 1758)'
$ws.Range("Q3").Value = 442085
$ws.Range("R3").Value = 7039138
$ws.Range("AC3").Value = 'Ringhack'

# Row 4
$ws.Range("A4").Value = 131064773
$ws.Range("Q4").Value = 442108
$ws.Range("AC4").Value = 'Ringhack äldre'

# Row 5
$ws.Range("A5").Value = 131064776
$ws.Range("Q5").Value = 442082
$ws.Range("R5").Value = 7039140
$ws.Range("AC5").Value = 'Ringhack färska och något äldre'

# Row 6
$ws.Range("A6").Value = 131064768
$ws.Range("Q6").Value = 442209
$ws.Range("R6").Value = 7039151
$ws.Range("AC6").Value = 'Ringhack färska och äldre'

# Row 9
$ws.Range("A9").Value = 131064778
$ws.Range("Q9").Value = 442145
$ws.Range("R9").Value = 7039101
$ws.Range("AC9").Value = 'Ringhack'

# Row 10
$ws.Range("A10").Value = 131064770
$ws.Range("Q10").Value = 442198
$ws.Range("R10").Value = 7039206
$ws.Range("AC10").Value = 'Ringhack färska och äldre'

# Row 11
$ws.Range("A11").Value = 131064783
$ws.Range("B11").Value = 91828
$ws.Range("E11").Value = 5432
$ws.Range("F11").Value = 'Granticka'
$ws.Range("G11").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H11").ClearContents()
$ws.Range("Q11").Value = 442292
$ws.Range("R11").Value = 7039182
$ws.Range("AC11").ClearContents()

# Row 12
$ws.Range("A12").Value = 131064766
$ws.Range("B12").Value = 57884
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = 'Tretåig hackspett'
$ws.Range("G12").Value = 'Picoides tridactylus'
$ws.Range("H12").Value = '(Linnaeus, 1758)'
$ws.Range("Q12").Value = 442271
$ws.Range("R12").Value = 7039174
$ws.Range("AC12").Value = 'Ringhack äldre'

# Row 13
$ws.Range("A13").Value = 131064763
$ws.Range("B13").Value = 57884
$ws.Range("E13").Value = 100109
$ws.Range("F13").Value = 'Tretåig hackspett'
$ws.Range("G13").Value = 'Picoides tridactylus'
$ws.Range("H13").Value = '(Linnaeus, 1758)'
$ws.Range("Q13").Value = 442230
$ws.Range("R13").Value = 7039147
$ws.Range("AC13").Value = 'Ringhack äldre'

# Row 14
$ws.Range("A14").Value = 131064779
$ws.Range("B14").Value = 91804
$ws.Range("E14").Value = 1108
$ws.Range("F14").Value = 'Harticka'
$ws.Range("G14").Value = 'Pelloporus leporinus'
$ws.Range("H14").Value = '(Fr.) Krieglst.'
$ws.Range("Q14").Value = 442245
$ws.Range("R14").Value = 7039149
$ws.Range("AC14").ClearContents()

# Row 18
$ws.Range("A18").Value = 131064781
$ws.Range("Q18").Value = 442200
$ws.Range("R18").Value = 7039150

# Row 19
$ws.Range("A19").Value = 131064780
$ws.Range("Q19").Value = 442259
$ws.Range("R19").Value = 7039181
